# Update Rizka - Update Data Binding Cari Mobil Baru 2
#
# The "Baru" sheet's first data row (A2) was bound to "Toyota Rush 2019";
# change it to just "Toyota Rush" (matching the shorter naming convention
# used elsewhere, e.g. "Daihatsu Sigra"/"Honda Brio"). Also leave the
# workbook focused on the "Baru" sheet with A3 selected, as it was left
# after making/reviewing this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baru")

$ws.Range("A2").Value = "Toyota Rush"

$ws.Activate()
$ws.Range("A3").Select()
